# Update "想去人数" (F column) figures on both the "展览" and "全部类型"
# sheets to reflect the latest scrape (commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 12934   # 南宁·第十九届（2024）良牙动漫夏季盛典（良牙夏典）: 12918 -> 12934
$ws1.Range("F5").Value = 633     # 南宁·蔚蓝档案only: 631 -> 633
$ws1.Range("F7").Value = 401     # 南宁·熊喵M动漫嘉年华【免费】: 400 -> 401
$ws1.Range("F8").Value = 1221    # 南宁·第二届北极光动漫展: 1216 -> 1221

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 12934   # 南宁·第十九届（2024）良牙动漫夏季盛典（良牙夏典）: 12918 -> 12934
$ws4.Range("F6").Value = 633     # 南宁·蔚蓝档案only: 631 -> 633
$ws4.Range("F10").Value = 401    # 南宁·熊喵M动漫嘉年华【免费】: 400 -> 401
$ws4.Range("F11").Value = 1221   # 南宁·第二届北极光动漫展: 1216 -> 1221
